# Update the "Förändrad" (Changed) date column (C) for rows 2-20
# from 2023-09-05 (serial 45174) to 2023-09-06 (serial 45175).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 20; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value = 45175
    }
}
